# Refactor timetable generation to apply consistent cell alignment and
# borders, and shift the Saturday (column G) concert-related entries down
# by one row (rows 19-28) to correct their timeslots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-flow the column-G content for the concert block: every entry from
#    row 19 down to row 28 moves one row later. Capture the text we need
#    before touching merges / values so nothing is lost along the way.
# ---------------------------------------------------------------------
$concertCallTime = $ws.Range("G19").Value2
$linaConcert     = $ws.Range("G20").Value2
$afterConcert    = $ws.Range("G26").Value2

# Break apart the merged ranges that are being resized so the cells
# inside them are free to be written individually.
$ws.Range("G15:G18").UnMerge()
$ws.Range("G20:G25").UnMerge()
$ws.Range("G26:G28").UnMerge()

# Clear the old homes of the text that is moving ...
$ws.Range("G19").Value2 = $null
$ws.Range("G26").Value2 = $null

# ... and write it into its new home, one row down.
$ws.Range("G20").Value2 = $concertCallTime
$ws.Range("G21").Value2 = $linaConcert
$ws.Range("G27").Value2 = $afterConcert

# Re-merge the (now shifted / resized) blocks.
$ws.Range("G15:G19").Merge()
$ws.Range("G21:G26").Merge()
$ws.Range("G27:G28").Merge()

# Excel's Merge() redistributes the border around merged ranges (only the
# outer edge stays "thin", inner seams get dropped). Reassert a uniform
# thin border on every side of every cell in the ranges we just
# re-merged so they keep looking the same as the rest of the table.
$xlContinuous = 1
$xlThin = 2
foreach ($ref in @("G15:G19", "G21:G26", "G27:G28")) {
    $rng = $ws.Range($ref)
    $rng.Borders.LineStyle = $xlContinuous
    $rng.Borders.Weight = $xlThin
}

# ---------------------------------------------------------------------
# 2) Apply consistent formatting across the whole table: every cell gets
#    centered horizontal + vertical alignment with text wrapping (the
#    border is already present on every cell of the used range). This
#    consolidates the old mix of "border only", "border + vertical
#    center + wrap" and "border + wrap" looks into a single uniform
#    style, same as is applied to the bold header row.
# ---------------------------------------------------------------------
$xlCenter = -4108
$body = $ws.Range("A1:G30")
$body.HorizontalAlignment = $xlCenter
$body.VerticalAlignment = $xlCenter
$body.WrapText = $true
